$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 2.03
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
